$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Athelete Id"
$ws.Range("B3").Value = "Ronda Rousey"
$ws.Range("B9").Value = "5'7"""
$ws.Range("B11").Value = "Glendale Fighting Club"
$ws.Range("B8").Value = "Women's Bantamweight"
$ws.Range("B14").Value = "USA"
$ws.Range("B17").Value = "Rowdy"
$ws.Range("B18").Value = "Orthodox"
$ws.Range("B20").Value = "Judo"
$ws.Range("B23").Value = "UFC"
$ws.Range("B24").Value = "UFClogo.jpg"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B32").Value = "2-0"
$ws.Range("B33").Value = "9-0"
$ws.Range("B15").Value = "2/1/1987"
$ws.Range("B6").Value = "United States"

$ws.Range("B10").Value = 135
$ws.Range("B16").Value = 28
$ws.Range("B19").Value = 66
$ws.Range("B29").Value = 11
$ws.Range("B30").Value = 0
$ws.Range("B31").Value = 0

$ws.Range("B7").Select()
